$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(61, 8).Value = 19818.334  # H61 was 1950
$ws.Cells.Item(61, 9).Value = 19818.334  # I61 was 1950
$ws.Cells.Item(61, 11).Value = 59455.00199999999  # K61 was 5850
$ws.Cells.Item(61, 13).Value = -59283.00199999999  # M61 was -5678
$ws.Cells.Item(113, 8).Value = 7752.154  # H113 was 7898.3335
$ws.Cells.Item(113, 10).Value = 7919.8887  # J113 was 8160.125
$ws.Cells.Item(113, 12).Value = 7919.8887  # L113 was 8160.125
$ws.Cells.Item(113, 14).Value = -14427.8887  # N113 was -14668.125
$ws.Cells.Item(116, 8).Value = 4685.125  # H116 was 5211.5713
$ws.Cells.Item(116, 9).Value = 3331.6667  # I116 was 4497.5
$ws.Cells.Item(116, 11).Value = 3331.6667  # K116 was 4497.5
$ws.Cells.Item(116, 13).Value = 110.3332999999998  # M116 was -1055.5
$ws.Cells.Item(137, 8).Value = 4141.3125  # H137 was 2984.4092
$ws.Cells.Item(137, 9).Value = 3788.4285  # I137 was 2171.9285
$ws.Cells.Item(137, 10).Value = 4415.778  # J137 was 4406.25
$ws.Cells.Item(137, 11).Value = 11365.2855  # K137 was 6515.7855
$ws.Cells.Item(137, 12).Value = 13247.334  # L137 was 13218.75
$ws.Cells.Item(137, 13).Value = -8815.2855  # M137 was -3965.7855
$ws.Cells.Item(137, 14).Value = -18347.334  # N137 was -18318.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4124.8423  # H32 was 4275.3335
$ws.Cells.Item(32, 10).Value = 33332.668  # J32 was 27498
$ws.Cells.Item(32, 12).Value = 33332.668  # L32 was 27498
$ws.Cells.Item(32, 14).Value = -33906.668  # N32 was -28072
$ws.Cells.Item(74, 8).Value = 2333.0667  # H74 was 2537.3076
$ws.Cells.Item(74, 9).Value = 2000.1666  # I74 was 2090
$ws.Cells.Item(74, 10).Value = 3664.6667  # J74 was 4997.5
$ws.Cells.Item(74, 11).Value = 2000.1666  # K74 was 2090
$ws.Cells.Item(74, 12).Value = 3664.6667  # L74 was 4997.5
$ws.Cells.Item(74, 13).Value = -1126.1666  # M74 was -1216
$ws.Cells.Item(74, 14).Value = -5412.6667  # N74 was -6745.5
$ws.Cells.Item(77, 8).Value = 2333.0667  # H77 was 2537.3076
$ws.Cells.Item(77, 9).Value = 2000.1666  # I77 was 2090
$ws.Cells.Item(77, 10).Value = 3664.6667  # J77 was 4997.5
$ws.Cells.Item(77, 11).Value = 10000.833  # K77 was 10450
$ws.Cells.Item(77, 12).Value = 18323.3335  # L77 was 24987.5
$ws.Cells.Item(77, 13).Value = -5632.833000000001  # M77 was -6082
$ws.Cells.Item(77, 14).Value = -27059.3335  # N77 was -33723.5
$ws.Cells.Item(110, 8).Value = 2372.5  # H110 was 2599.1667
$ws.Cells.Item(110, 9).Value = 2330  # I110 was 2648.75
$ws.Cells.Item(110, 11).Value = 2330  # K110 was 2648.75
$ws.Cells.Item(110, 13).Value = -285  # M110 was -603.75
$ws.Cells.Item(122, 8).Value = 2898.25  # H122 was 2926.8572
$ws.Cells.Item(122, 9).Value = 2898.25  # I122 was 2926.8572
$ws.Cells.Item(122, 11).Value = 8694.75  # K122 was 8780.571599999999
$ws.Cells.Item(122, 13).Value = -6244.75  # M122 was -6330.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 750.3333  # H11 was 600.4
$ws.Cells.Item(11, 10).Value = 1749  # J11 was 1998
$ws.Cells.Item(11, 12).Value = 1749  # L11 was 1998
$ws.Cells.Item(11, 14).Value = -2029  # N11 was -2278
$ws.Cells.Item(31, 8).Value = 15375  # H31 was 17265
$ws.Cells.Item(31, 10).Value = 20000  # J31 was 20522.5
$ws.Cells.Item(31, 12).Value = 20000  # L31 was 20522.5
$ws.Cells.Item(31, 14).Value = -20504  # N31 was -21026.5
$ws.Cells.Item(107, 8).Value = 865  # H107 was 684.8
$ws.Cells.Item(107, 9).Value = 872.5  # I107 was 684.8
$ws.Cells.Item(107, 10).Value = 850  # J107 was 0
$ws.Cells.Item(107, 11).Value = 872.5  # K107 was 684.8
$ws.Cells.Item(107, 12).Value = 850  # L107 was 0
$ws.Cells.Item(107, 13).Value = 1047.5  # M107 was 1235.2
$ws.Cells.Item(107, 14).Value = -4690  # N107 was None

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0  # H4 was 1000
$ws.Cells.Item(4, 10).Value = 0  # J4 was 1000
$ws.Cells.Item(4, 12).Value = 0  # L4 was 1000
$ws.Cells.Item(4, 14).ClearContents()  # N4 was -1224
$ws.Cells.Item(16, 8).Value = 1750  # H16 was 551.1111
$ws.Cells.Item(16, 9).Value = 1750  # I16 was 475.33334
$ws.Cells.Item(16, 10).Value = 0  # J16 was 702.6667
$ws.Cells.Item(16, 11).Value = 1750  # K16 was 475.33334
$ws.Cells.Item(16, 12).Value = 0  # L16 was 702.6667
$ws.Cells.Item(16, 13).Value = -1463  # M16 was -188.33334
$ws.Cells.Item(16, 14).ClearContents()  # N16 was -1276.6667
$ws.Cells.Item(31, 8).Value = 2073.4583  # H31 was 2077.625
$ws.Cells.Item(31, 10).Value = 2672.7  # J31 was 2682.7
$ws.Cells.Item(31, 12).Value = 2672.7  # L31 was 2682.7
$ws.Cells.Item(31, 14).Value = -3262.7  # N31 was -3272.7
$ws.Cells.Item(34, 8).Value = 2073.4583  # H34 was 2077.625
$ws.Cells.Item(34, 10).Value = 2672.7  # J34 was 2682.7
$ws.Cells.Item(34, 12).Value = 2672.7  # L34 was 2682.7
$ws.Cells.Item(34, 14).Value = -3076.7  # N34 was -3086.7
$ws.Cells.Item(41, 8).Value = 10575.357  # H41 was 9005.182000000001
$ws.Cells.Item(41, 10).Value = 12083  # J41 was 10666.444
$ws.Cells.Item(41, 12).Value = 12083  # L41 was 10666.444
$ws.Cells.Item(41, 14).Value = -12939  # N41 was -11522.444
$ws.Cells.Item(99, 8).Value = 1999.75  # H99 was 1570.8572
$ws.Cells.Item(99, 9).Value = 1833  # I99 was 1499.2
$ws.Cells.Item(99, 10).Value = 2500  # J99 was 1750
$ws.Cells.Item(99, 11).Value = 1833  # K99 was 1499.2
$ws.Cells.Item(99, 12).Value = 2500  # L99 was 1750
$ws.Cells.Item(99, 13).Value = -335  # M99 was -1.200000000000045
$ws.Cells.Item(99, 14).Value = -5496  # N99 was -4746
$ws.Cells.Item(107, 8).Value = 676.1875  # H107 was 665.3125
$ws.Cells.Item(107, 9).Value = 663.0769  # I107 was 638.9286
$ws.Cells.Item(107, 10).Value = 733  # J107 was 850
$ws.Cells.Item(107, 11).Value = 663.0769  # K107 was 638.9286
$ws.Cells.Item(107, 12).Value = 733  # L107 was 850
$ws.Cells.Item(107, 13).Value = 1256.9231  # M107 was 1281.0714
$ws.Cells.Item(107, 14).Value = -4573  # N107 was -4690
$ws.Cells.Item(113, 8).Value = 1750  # H113 was 551.1111
$ws.Cells.Item(113, 9).Value = 1750  # I113 was 475.33334
$ws.Cells.Item(113, 10).Value = 0  # J113 was 702.6667
$ws.Cells.Item(113, 11).Value = 1750  # K113 was 475.33334
$ws.Cells.Item(113, 12).Value = 0  # L113 was 702.6667
$ws.Cells.Item(113, 13).Value = 420  # M113 was 1694.66666
$ws.Cells.Item(113, 14).ClearContents()  # N113 was -5042.6667
$ws.Cells.Item(126, 8).Value = 1999.75  # H126 was 1570.8572
$ws.Cells.Item(126, 9).Value = 1833  # I126 was 1499.2
$ws.Cells.Item(126, 10).Value = 2500  # J126 was 1750
$ws.Cells.Item(126, 11).Value = 5499  # K126 was 4497.6
$ws.Cells.Item(126, 12).Value = 7500  # L126 was 5250
$ws.Cells.Item(126, 13).Value = -3029  # M126 was -2027.6
$ws.Cells.Item(126, 14).Value = -12440  # N126 was -10190

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 703.75  # H2 was 900
$ws.Cells.Item(2, 9).Value = 407.5  # I2 was 800
$ws.Cells.Item(2, 11).Value = 2445  # K2 was 4800
$ws.Cells.Item(2, 13).Value = -2332  # M2 was -4687
$ws.Cells.Item(23, 8).Value = 779.4286  # H23 was 900.8333
$ws.Cells.Item(23, 9).Value = 384  # I23 was 550.5
$ws.Cells.Item(23, 11).Value = 1152  # K23 was 1651.5
$ws.Cells.Item(23, 13).Value = -917  # M23 was -1416.5
$ws.Cells.Item(64, 8).Value = 500  # H64 was 0
$ws.Cells.Item(64, 9).Value = 500  # I64 was 0
$ws.Cells.Item(64, 11).Value = 1500  # K64 was 0
$ws.Cells.Item(64, 13).Value = -1230  # M64 was None
$ws.Cells.Item(67, 8).Value = 500  # H67 was 0
$ws.Cells.Item(67, 9).Value = 500  # I67 was 0
$ws.Cells.Item(67, 11).Value = 1500  # K67 was 0
$ws.Cells.Item(67, 13).Value = -564  # M67 was None
$ws.Cells.Item(93, 8).Value = 0  # H93 was 327
$ws.Cells.Item(93, 10).Value = 0  # J93 was 327
$ws.Cells.Item(93, 12).Value = 0  # L93 was 981
$ws.Cells.Item(93, 14).ClearContents()  # N93 was -4725
$ws.Cells.Item(140, 8).Value = 2564.2856  # H140 was 2331.4443
$ws.Cells.Item(140, 9).Value = 1992.1666  # I140 was 1873.25
$ws.Cells.Item(140, 11).Value = 5976.4998  # K140 was 5619.75
$ws.Cells.Item(140, 13).Value = -796.4997999999996  # M140 was -439.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 0  # H35 was 3025
$ws.Cells.Item(35, 9).Value = 0  # I35 was 3025
$ws.Cells.Item(35, 11).Value = 0  # K35 was 3025
$ws.Cells.Item(35, 13).ClearContents()  # M35 was -2727
$ws.Cells.Item(41, 8).Value = 6116.6665  # H41 was 3362.75
$ws.Cells.Item(41, 9).Value = 4800  # I41 was 3450.3333
$ws.Cells.Item(41, 10).Value = 8750  # J41 was 3100
$ws.Cells.Item(41, 11).Value = 4800  # K41 was 3450.3333
$ws.Cells.Item(41, 12).Value = 8750  # L41 was 3100
$ws.Cells.Item(41, 13).Value = -4445  # M41 was -3095.3333
$ws.Cells.Item(41, 14).Value = -9460  # N41 was -3810
$ws.Cells.Item(80, 8).Value = 4000  # H80 was 950
$ws.Cells.Item(80, 9).Value = 0  # I80 was 950
$ws.Cells.Item(80, 10).Value = 4000  # J80 was 0
$ws.Cells.Item(80, 11).Value = 0  # K80 was 950
$ws.Cells.Item(80, 12).Value = 4000  # L80 was 0
$ws.Cells.Item(80, 13).ClearContents()  # M80 was 48
$ws.Cells.Item(80, 14).Value = -5996  # N80 was None
$ws.Cells.Item(83, 8).Value = 4000  # H83 was 950
$ws.Cells.Item(83, 9).Value = 0  # I83 was 950
$ws.Cells.Item(83, 10).Value = 4000  # J83 was 0
$ws.Cells.Item(83, 11).Value = 0  # K83 was 4750
$ws.Cells.Item(83, 12).Value = 20000  # L83 was 0
$ws.Cells.Item(83, 13).ClearContents()  # M83 was 242
$ws.Cells.Item(83, 14).Value = -29984  # N83 was None
$ws.Cells.Item(113, 8).Value = 0  # H113 was 1075
$ws.Cells.Item(113, 9).Value = 0  # I113 was 1075
$ws.Cells.Item(113, 11).Value = 0  # K113 was 1075
$ws.Cells.Item(113, 13).ClearContents()  # M113 was 1095
$ws.Cells.Item(122, 8).Value = 8613.625  # H122 was 8738.8125
$ws.Cells.Item(122, 9).Value = 10344.417  # I122 was 10511.333
$ws.Cells.Item(122, 11).Value = 31033.251  # K122 was 31533.999
$ws.Cells.Item(122, 13).Value = -28583.251  # M122 was -29083.999
$ws.Cells.Item(126, 8).Value = 8372.362999999999  # H126 was 8372.454
$ws.Cells.Item(126, 9).Value = 5579.2  # I126 was 6624.25
$ws.Cells.Item(126, 10).Value = 10700  # J126 was 9371.429
$ws.Cells.Item(126, 11).Value = 16737.6  # K126 was 19872.75
$ws.Cells.Item(126, 12).Value = 32100  # L126 was 28114.287
$ws.Cells.Item(126, 13).Value = -14267.6  # M126 was -17402.75
$ws.Cells.Item(126, 14).Value = -37040  # N126 was -33054.287
$ws.Cells.Item(132, 8).Value = 866.3333  # H132 was 900
$ws.Cells.Item(132, 9).Value = 866.3333  # I132 was 900
$ws.Cells.Item(132, 11).Value = 2598.9999  # K132 was 2700
$ws.Cells.Item(132, 13).Value = -68.9998999999998  # M132 was -170

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1150  # H22 was 1072.6666
$ws.Cells.Item(22, 9).Value = 1000  # I22 was 921.5
$ws.Cells.Item(22, 11).Value = 1000  # K22 was 921.5
$ws.Cells.Item(22, 13).Value = -705  # M22 was -626.5
$ws.Cells.Item(27, 8).Value = 1150  # H27 was 1072.6666
$ws.Cells.Item(27, 9).Value = 1000  # I27 was 921.5
$ws.Cells.Item(27, 11).Value = 1000  # K27 was 921.5
$ws.Cells.Item(27, 13).Value = -893  # M27 was -814.5
$ws.Cells.Item(40, 8).Value = 2734.6667  # H40 was 4650
$ws.Cells.Item(40, 9).Value = 1602  # I40 was 1800
$ws.Cells.Item(40, 10).Value = 5000  # J40 was 7500
$ws.Cells.Item(40, 11).Value = 1602  # K40 was 1800
$ws.Cells.Item(40, 12).Value = 5000  # L40 was 7500
$ws.Cells.Item(40, 13).Value = -1466  # M40 was -1664
$ws.Cells.Item(40, 14).Value = -5272  # N40 was -7772
$ws.Cells.Item(61, 8).Value = 0  # H61 was 3999
$ws.Cells.Item(61, 9).Value = 0  # I61 was 3999
$ws.Cells.Item(61, 11).Value = 0  # K61 was 3999
$ws.Cells.Item(61, 13).ClearContents()  # M61 was -3797
$ws.Cells.Item(68, 8).Value = 46785.145  # H68 was 46856.57
$ws.Cells.Item(68, 10).Value = 79373.5  # J68 was 79498.5
$ws.Cells.Item(68, 12).Value = 79373.5  # L68 was 79498.5
$ws.Cells.Item(68, 14).Value = -80871.5  # N68 was -80996.5
$ws.Cells.Item(71, 8).Value = 46785.145  # H71 was 46856.57
$ws.Cells.Item(71, 10).Value = 79373.5  # J71 was 79498.5
$ws.Cells.Item(71, 12).Value = 396867.5  # L71 was 397492.5
$ws.Cells.Item(71, 14).Value = -404355.5  # N71 was -404980.5
$ws.Cells.Item(109, 8).Value = 0  # H109 was 69990
$ws.Cells.Item(109, 10).Value = 0  # J109 was 69990
$ws.Cells.Item(109, 12).Value = 0  # L109 was 69990
$ws.Cells.Item(109, 14).ClearContents()  # N109 was -72764
$ws.Cells.Item(113, 8).Value = 0  # H113 was 3999
$ws.Cells.Item(113, 9).Value = 0  # I113 was 3999
$ws.Cells.Item(113, 11).Value = 0  # K113 was 3999
$ws.Cells.Item(113, 13).ClearContents()  # M113 was -1829

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 502.5  # H14 was 0
$ws.Cells.Item(14, 9).Value = 500  # I14 was 0
$ws.Cells.Item(14, 10).Value = 505  # J14 was 0
$ws.Cells.Item(14, 11).Value = 500  # K14 was 0
$ws.Cells.Item(14, 12).Value = 505  # L14 was 0
$ws.Cells.Item(14, 13).Value = -332  # M14 was None
$ws.Cells.Item(14, 14).Value = -841  # N14 was None
$ws.Cells.Item(80, 8).Value = 0  # H80 was 21425
$ws.Cells.Item(80, 10).Value = 0  # J80 was 21425
$ws.Cells.Item(80, 12).Value = 0  # L80 was 21425
$ws.Cells.Item(80, 14).ClearContents()  # N80 was -23421
$ws.Cells.Item(83, 8).Value = 0  # H83 was 21425
$ws.Cells.Item(83, 10).Value = 0  # J83 was 21425
$ws.Cells.Item(83, 12).Value = 0  # L83 was 64275
$ws.Cells.Item(83, 14).ClearContents()  # N83 was -74259
$ws.Cells.Item(107, 8).Value = 1178.0667  # H107 was 1179.6666
$ws.Cells.Item(107, 9).Value = 842.8182  # I107 was 859.6
$ws.Cells.Item(107, 10).Value = 2100  # J107 was 1819.8
$ws.Cells.Item(107, 11).Value = 2528.4546  # K107 was 2578.8
$ws.Cells.Item(107, 12).Value = 6300  # L107 was 5459.4
$ws.Cells.Item(107, 13).Value = -608.4546  # M107 was -658.8000000000002
$ws.Cells.Item(107, 14).Value = -10140  # N107 was -9299.4
$ws.Cells.Item(113, 8).Value = 0  # H113 was 269
$ws.Cells.Item(113, 9).Value = 0  # I113 was 269
$ws.Cells.Item(113, 11).Value = 0  # K113 was 807
$ws.Cells.Item(113, 13).ClearContents()  # M113 was 1363
$ws.Cells.Item(122, 8).Value = 2833  # H122 was 3000
$ws.Cells.Item(122, 9).Value = 2999  # I122 was 0
$ws.Cells.Item(122, 10).Value = 2750  # J122 was 3000
$ws.Cells.Item(122, 11).Value = 8997  # K122 was 0
$ws.Cells.Item(122, 12).Value = 8250  # L122 was 9000
$ws.Cells.Item(122, 13).Value = -6547  # M122 was None
$ws.Cells.Item(122, 14).Value = -13150  # N122 was -13900
$ws.Cells.Item(132, 8).Value = 2951.7273  # H132 was 3683.75
$ws.Cells.Item(132, 9).Value = 1219.6  # I132 was 1397
$ws.Cells.Item(132, 10).Value = 4395.1665  # J132 was 4446
$ws.Cells.Item(132, 11).Value = 3658.8  # K132 was 4191
$ws.Cells.Item(132, 12).Value = 13185.4995  # L132 was 13338
$ws.Cells.Item(132, 13).Value = -1128.8  # M132 was -1661
$ws.Cells.Item(132, 14).Value = -18245.4995  # N132 was -18398
$ws.Cells.Item(136, 8).Value = 9472.416999999999  # H136 was 10464.111
$ws.Cells.Item(136, 9).Value = 9031.666999999999  # I136 was 10298.833
$ws.Cells.Item(136, 11).Value = 27095.001  # K136 was 30896.499
$ws.Cells.Item(136, 13).Value = -24545.001  # M136 was -28346.499
